$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Ano safra"
$ws.Range("D1").Value = "Estoque Inicial(mi)"

$ws.Range("J6").Select()
